# User checkpoint: Add customer data validation and importer method.
#
# The "Template" sheet currently has 3 header columns:
#   A1 customer_name [Data]
#   B1 customer_group [Link [Customer Group]]
#   C1 territory [Link [Territory]]
#
# Target layout (A1:H1):
#   A1 customer_name [Data]                                   (unchanged)
#   B1 customer_type [Select]                                 (new)
#   C1 customer_group [Link [Customer Group]]                 (was B1)
#   D1 territory [Link [Territory]]                           (was C1)
#   E1 default_currency [Link [Currency]]                     (new)
#   F1 default_price_list [Link [Price List]]                 (new)
#   G1 tax_id [Data]                                           (new)
#   H1 payment_terms [Link [Payment Terms Template]]           (new)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one new column before column B. This shifts the existing
# customer_group/territory headers from B/C to C/D, and - importantly -
# the inserted column inherits the header cell's style (bold, centered,
# bordered) from the columns around it, matching the existing header look.
$ws.Columns("B").Insert()
$ws.Range("B1").Value = "customer_type [Select]"

# Columns E:H are brand new (beyond the old used range), so their cells
# don't automatically pick up the header style. Copy the style from the
# existing A1 header cell onto E1:H1 before filling in the values.
$ws.Range("A1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)

$ws.Range("E1").Value = "default_currency [Link [Currency]]"
$ws.Range("F1").Value = "default_price_list [Link [Price List]]"
$ws.Range("G1").Value = "tax_id [Data]"
$ws.Range("H1").Value = "payment_terms [Link [Payment Terms Template]]"
